# Applies the categorized-trial regeneration described in the commit:
# input list was regenerated (new subject-level trial numbering, new stimuli
# draws, refreshed rating/n/percentile stats) for rows 2-27 of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
  @{r=2; c=3; v=3},
  @{r=2; c=6; v=109},
  @{r=2; c=12; v="stimuli/img_0kqc0.png"},
  @{r=2; c=13; v=43.74468085106383},
  @{r=2; c=14; v=27.14893617021277},
  @{r=2; c=15; v=35.4468085106383},
  @{r=2; c=16; v=47},
  @{r=2; c=17; v=2},
  @{r=2; c=18; v=2},
  @{r=2; c=19; v=2},
  @{r=3; c=3; v=3},
  @{r=3; c=6; v=110},
  @{r=3; c=12; v="stimuli/img_bbs77.png"},
  @{r=3; c=13; v=31.64444444444445},
  @{r=3; c=14; v=21.26666666666667},
  @{r=3; c=15; v=26.45555555555556},
  @{r=3; c=16; v=45},
  @{r=4; c=3; v=3},
  @{r=4; c=6; v=111},
  @{r=4; c=12; v="stimuli/img_w8yhd.png"},
  @{r=4; c=13; v=55.74418604651163},
  @{r=4; c=14; v=38.90697674418605},
  @{r=4; c=15; v=47.32558139534883},
  @{r=4; c=16; v=43},
  @{r=4; c=17; v=4},
  @{r=4; c=18; v=4},
  @{r=4; c=19; v=4},
  @{r=5; c=3; v=3},
  @{r=5; c=6; v=112},
  @{r=5; c=12; v="stimuli/img_tujn3.png"},
  @{r=5; c=13; v=81.40909090909091},
  @{r=5; c=14; v=62.52272727272727},
  @{r=5; c=15; v=71.96590909090909},
  @{r=5; c=16; v=44},
  @{r=5; c=17; v=8},
  @{r=5; c=18; v=8},
  @{r=5; c=19; v=8},
  @{r=6; c=3; v=3},
  @{r=6; c=6; v=113},
  @{r=6; c=12; v="stimuli/img_6a0hu.png"},
  @{r=6; c=13; v=61.275},
  @{r=6; c=14; v=42.025},
  @{r=6; c=15; v=51.65},
  @{r=6; c=16; v=40},
  @{r=6; c=17; v=4},
  @{r=6; c=18; v=4},
  @{r=6; c=19; v=4},
  @{r=7; c=3; v=3},
  @{r=7; c=6; v=114},
  @{r=7; c=8; v="kitchens"},
  @{r=7; c=12; v="stimuli/img_7wquy.png"},
  @{r=7; c=13; v=50.59375},
  @{r=7; c=14; v=30.59375},
  @{r=7; c=15; v=40.59375},
  @{r=7; c=16; v=32},
  @{r=7; c=17; v=2},
  @{r=7; c=18; v=2},
  @{r=7; c=19; v=2},
  @{r=8; c=3; v=3},
  @{r=8; c=6; v=115},
  @{r=8; c=8; v="living_rooms"},
  @{r=8; c=9; v="target"},
  @{r=8; c=11; v="j"},
  @{r=8; c=12; v="stimuli/img_6zz63.png"},
  @{r=8; c=13; v=87.66666666666667},
  @{r=8; c=14; v=70.59999999999999},
  @{r=8; c=15; v=79.13333333333333},
  @{r=8; c=16; v=45},
  @{r=8; c=17; v=9},
  @{r=8; c=18; v=10},
  @{r=8; c=19; v=10},
  @{r=9; c=3; v=3},
  @{r=9; c=6; v=116},
  @{r=9; c=12; v="stimuli/img_xy930.png"},
  @{r=9; c=13; v=70.5952380952381},
  @{r=9; c=14; v=49.47619047619047},
  @{r=9; c=15; v=60.03571428571429},
  @{r=9; c=16; v=42},
  @{r=9; c=17; v=6},
  @{r=9; c=18; v=6},
  @{r=9; c=19; v=6},
  @{r=10; c=3; v=3},
  @{r=10; c=6; v=117},
  @{r=10; c=8; v="bedrooms"},
  @{r=10; c=9; v="distractor"},
  @{r=10; c=11; v="f"},
  @{r=10; c=12; v="stimuli/img_die1d.png"},
  @{r=10; c=13; v=75.42857142857143},
  @{r=10; c=14; v=53.30952380952381},
  @{r=10; c=15; v=64.36904761904762},
  @{r=10; c=16; v=42},
  @{r=10; c=17; v=6},
  @{r=10; c=18; v=6},
  @{r=10; c=19; v=6},
  @{r=11; c=3; v=3},
  @{r=11; c=6; v=118},
  @{r=11; c=8; v="kitchens"},
  @{r=11; c=9; v="distractor"},
  @{r=11; c=11; v="f"},
  @{r=11; c=12; v="stimuli/img_gztbt.png"},
  @{r=11; c=13; v=55.06451612903226},
  @{r=11; c=14; v=26.09677419354839},
  @{r=11; c=15; v=40.58064516129032},
  @{r=11; c=16; v=31},
  @{r=11; c=17; v=2},
  @{r=11; c=18; v=2},
  @{r=11; c=19; v=2},
  @{r=12; c=3; v=3},
  @{r=12; c=6; v=119},
  @{r=12; c=12; v="stimuli/img_kost0.png"},
  @{r=12; c=13; v=63.09090909090909},
  @{r=12; c=14; v=42.77272727272727},
  @{r=12; c=15; v=52.93181818181819},
  @{r=12; c=16; v=44},
  @{r=12; c=17; v=5},
  @{r=12; c=18; v=5},
  @{r=12; c=19; v=5},
  @{r=13; c=3; v=3},
  @{r=13; c=6; v=120},
  @{r=13; c=8; v="kitchens"},
  @{r=13; c=9; v="distractor"},
  @{r=13; c=11; v="f"},
  @{r=13; c=12; v="stimuli/img_pt3d7.png"},
  @{r=13; c=13; v=65.08571428571429},
  @{r=13; c=14; v=44.65714285714286},
  @{r=13; c=15; v=54.87142857142857},
  @{r=13; c=16; v=35},
  @{r=13; c=17; v=4},
  @{r=13; c=18; v=4},
  @{r=13; c=19; v=4},
  @{r=14; c=3; v=3},
  @{r=14; c=6; v=121},
  @{r=14; c=8; v="living_rooms"},
  @{r=14; c=9; v="target"},
  @{r=14; c=11; v="j"},
  @{r=14; c=12; v="stimuli/img_95hiv.png"},
  @{r=14; c=13; v=84.04545454545455},
  @{r=14; c=14; v=67.31818181818181},
  @{r=14; c=15; v=75.68181818181819},
  @{r=14; c=16; v=44},
  @{r=14; c=17; v=9},
  @{r=14; c=18; v=9},
  @{r=14; c=19; v=9},
  @{r=15; c=3; v=3},
  @{r=15; c=6; v=122},
  @{r=15; c=8; v="bedrooms"},
  @{r=15; c=9; v="distractor"},
  @{r=15; c=11; v="f"},
  @{r=15; c=12; v="stimuli/img_5p2ql.png"},
  @{r=15; c=13; v=89.19565217391305},
  @{r=15; c=14; v=72.52173913043478},
  @{r=15; c=15; v=80.85869565217391},
  @{r=15; c=16; v=46},
  @{r=15; c=17; v=10},
  @{r=15; c=18; v=10},
  @{r=15; c=19; v=10},
  @{r=16; c=3; v=3},
  @{r=16; c=6; v=123},
  @{r=16; c=12; v="stimuli/img_di6f0.png"},
  @{r=16; c=13; v=94.04347826086956},
  @{r=16; c=14; v=83.34782608695652},
  @{r=16; c=15; v=88.69565217391303},
  @{r=16; c=16; v=46},
  @{r=16; c=17; v=10},
  @{r=16; c=18; v=10},
  @{r=16; c=19; v=10},
  @{r=17; c=3; v=3},
  @{r=17; c=6; v=124},
  @{r=17; c=12; v="stimuli/img_xbtev.png"},
  @{r=17; c=13; v=13.68181818181818},
  @{r=17; c=14; v=8.568181818181818},
  @{r=17; c=15; v=11.125},
  @{r=17; c=16; v=44},
  @{r=17; c=17; v=1},
  @{r=17; c=18; v=1},
  @{r=17; c=19; v=1},
  @{r=18; c=3; v=3},
  @{r=18; c=6; v=125},
  @{r=18; c=12; v="stimuli/img_pey7u.png"},
  @{r=18; c=13; v=30.34883720930232},
  @{r=18; c=14; v=20.34883720930232},
  @{r=18; c=15; v=25.34883720930232},
  @{r=18; c=16; v=43},
  @{r=18; c=17; v=1},
  @{r=18; c=18; v=2},
  @{r=18; c=19; v=2},
  @{r=19; c=3; v=3},
  @{r=19; c=6; v=126},
  @{r=19; c=12; v="stimuli/img_wz6x5.png"},
  @{r=19; c=13; v=68.3695652173913},
  @{r=19; c=14; v=48.47826086956522},
  @{r=19; c=15; v=58.42391304347826},
  @{r=19; c=16; v=46},
  @{r=19; c=17; v=5},
  @{r=19; c=18; v=5},
  @{r=19; c=19; v=5},
  @{r=20; c=3; v=3},
  @{r=20; c=6; v=127},
  @{r=20; c=8; v="living_rooms"},
  @{r=20; c=9; v="target"},
  @{r=20; c=11; v="j"},
  @{r=20; c=12; v="stimuli/img_xu1p3.png"},
  @{r=20; c=13; v=75.27659574468085},
  @{r=20; c=14; v=56.68085106382978},
  @{r=20; c=15; v=65.97872340425532},
  @{r=20; c=16; v=47},
  @{r=20; c=17; v=7},
  @{r=20; c=18; v=7},
  @{r=20; c=19; v=7},
  @{r=21; c=3; v=3},
  @{r=21; c=6; v=128},
  @{r=21; c=12; v="stimuli/img_wgkqa.png"},
  @{r=21; c=13; v=87.25581395348837},
  @{r=21; c=14; v=71.13953488372093},
  @{r=21; c=15; v=79.19767441860465},
  @{r=21; c=16; v=43},
  @{r=21; c=17; v=10},
  @{r=21; c=18; v=10},
  @{r=21; c=19; v=10},
  @{r=22; c=3; v=3},
  @{r=22; c=6; v=129},
  @{r=22; c=12; v="stimuli/img_abobq.png"},
  @{r=22; c=13; v=75.18421052631579},
  @{r=22; c=14; v=54.13157894736842},
  @{r=22; c=15; v=64.65789473684211},
  @{r=22; c=16; v=38},
  @{r=22; c=17; v=6},
  @{r=22; c=18; v=6},
  @{r=22; c=19; v=6},
  @{r=23; c=3; v=3},
  @{r=23; c=6; v=130},
  @{r=23; c=12; v="stimuli/img_eh0no.png"},
  @{r=23; c=13; v=53.66666666666666},
  @{r=23; c=14; v=36.02564102564103},
  @{r=23; c=15; v=44.84615384615385},
  @{r=23; c=16; v=39},
  @{r=23; c=17; v=3},
  @{r=23; c=18; v=3},
  @{r=23; c=19; v=3},
  @{r=24; c=3; v=3},
  @{r=24; c=6; v=131},
  @{r=24; c=12; v="stimuli/img_4o8l0.png"},
  @{r=24; c=13; v=46.02173913043478},
  @{r=24; c=14; v=31.45652173913043},
  @{r=24; c=15; v=38.73913043478261},
  @{r=24; c=16; v=46},
  @{r=24; c=17; v=3},
  @{r=24; c=18; v=3},
  @{r=24; c=19; v=3},
  @{r=25; c=3; v=3},
  @{r=25; c=6; v=132},
  @{r=25; c=12; v="stimuli/img_cehin.png"},
  @{r=25; c=13; v=78.86363636363636},
  @{r=25; c=14; v=60.02272727272727},
  @{r=25; c=15; v=69.44318181818181},
  @{r=25; c=16; v=44},
  @{r=25; c=17; v=7},
  @{r=25; c=18; v=7},
  @{r=25; c=19; v=7},
  @{r=26; c=3; v=3},
  @{r=26; c=6; v=133},
  @{r=26; c=12; v="stimuli/img_mdpr4.png"},
  @{r=26; c=13; v=74.04255319148936},
  @{r=26; c=14; v=54.70212765957447},
  @{r=26; c=15; v=64.37234042553192},
  @{r=26; c=16; v=47},
  @{r=26; c=17; v=6},
  @{r=26; c=18; v=6},
  @{r=26; c=19; v=6},
  @{r=27; c=3; v=3},
  @{r=27; c=6; v=134},
  @{r=27; c=8; v="living_rooms"},
  @{r=27; c=9; v="target"},
  @{r=27; c=11; v="j"},
  @{r=27; c=12; v="stimuli/img_bj99b.png"},
  @{r=27; c=13; v=82.79069767441861},
  @{r=27; c=14; v=65.46511627906976},
  @{r=27; c=15; v=74.12790697674419},
  @{r=27; c=16; v=43},
  @{r=27; c=17; v=8},
  @{r=27; c=18; v=8},
  @{r=27; c=19; v=8}
)

foreach ($item in $changes) {
  $ws.Cells.Item($item.r, $item.c).Value = $item.v
}
